# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update the "last updated" timestamp shown in the title cell (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 10 de Mayo de 2020 a las 05:34"

# --- Countries ranking shuffled: Pakistan overtakes Ecuador (rows 23-24) ---
$ws.Range("A23").Value = "Pakistan"
$ws.Range("B23").Value = 29465
$ws.Range("C23").Value = 729
$ws.Range("D23").Value = 8023
$ws.Range("E23").Value = 20803
$ws.Range("F23").Value = 111
$ws.Range("G23").Value = 3
$ws.Range("H23").Value = 639

$ws.Range("A24").Value = "Ecuador"
$ws.Range("B24").Value = 29071
$ws.Range("C24").Value = 0
$ws.Range("D24").Value = 3433
$ws.Range("E24").Value = 23921
$ws.Range("F24").Value = 181
$ws.Range("G24").Value = 0
$ws.Range("H24").Value = 1717

# --- El Salvador moves up above Albania/Sri Lanka/Niger/Libano/Maldivas (rows 102-107) ---
$ws.Range("A102").Value = "El Salvador"
$ws.Range("B102").Value = 889
$ws.Range("C102").Value = 105
$ws.Range("D102").Value = 293
$ws.Range("E102").Value = 579
$ws.Range("F102").Value = 4
$ws.Range("G102").Value = 0
$ws.Range("H102").Value = 17

$ws.Range("A103").Value = "Albania"
$ws.Range("B103").Value = 856
$ws.Range("C103").Value = 0
$ws.Range("D103").Value = 627
$ws.Range("E103").Value = 198
$ws.Range("F103").Value = 7
$ws.Range("G103").Value = 0
$ws.Range("H103").Value = 31

$ws.Range("A104").Value = "Sri Lanka"
$ws.Range("B104").Value = 847
$ws.Range("C104").Value = 0
$ws.Range("D104").Value = 260
$ws.Range("E104").Value = 578
$ws.Range("F104").Value = 1
$ws.Range("G104").Value = 0
$ws.Range("H104").Value = 9

$ws.Range("A105").Value = "Niger"
$ws.Range("B105").Value = 815
$ws.Range("C105").Value = 0
$ws.Range("D105").Value = 617
$ws.Range("E105").Value = 153
$ws.Range("F105").Value = 0
$ws.Range("G105").Value = 0
$ws.Range("H105").Value = 45

$ws.Range("A106").Value = "Libano"
$ws.Range("B106").Value = 809
$ws.Range("C106").Value = 0
$ws.Range("D106").Value = 234
$ws.Range("E106").Value = 549
$ws.Range("F106").Value = 4
$ws.Range("G106").Value = 0
$ws.Range("H106").Value = 26

$ws.Range("A107").Value = "Maldivas"
$ws.Range("B107").Value = 790
$ws.Range("C107").Value = 0
$ws.Range("D107").Value = 29
$ws.Range("E107").Value = 758
$ws.Range("F107").Value = 2
$ws.Range("G107").Value = 0
$ws.Range("H107").Value = 3

# --- Row 123 (Jamaica) updated counts, ranking unchanged ---
$ws.Range("B123").Value = 498
$ws.Range("C123").Value = 8
$ws.Range("D123").Value = 78
$ws.Range("E123").Value = 411
